$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C, rows 2 through 252 all currently hold 7573; change them to 7293.
$ws.Range("C2:C252").Value = 7293
